# Auto-generated: refresh cached market-data values in Odin_Profits sheets
# (scheduled runner update - plain numeric cache values, no formulas)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 15889.875
$ws.Range("I32").Value = 8989
$ws.Range("J32").Value = 16875.715
$ws.Range("K32").Value = 8989
$ws.Range("L32").Value = 16875.715
$ws.Range("M32").Value = -8663
$ws.Range("N32").Value = -17527.715
$ws.Range("H42").Value = 632
$ws.Range("I42").Value = 720.4286
$ws.Range("J42").Value = 13
$ws.Range("K42").Value = 2161.2858
$ws.Range("L42").Value = 39
$ws.Range("M42").Value = -1931.2858
$ws.Range("N42").Value = -499
$ws.Range("H76").Value = 41674796
$ws.Range("I76").Value = 66675476
$ws.Range("K76").Value = 66675476
$ws.Range("M76").Value = -66675161
$ws.Range("H79").Value = 41674796
$ws.Range("I79").Value = 66675476
$ws.Range("K79").Value = 66675476
$ws.Range("M79").Value = -66674384
$ws.Range("H97").Value = 2300
$ws.Range("J97").Value = 2300
$ws.Range("L97").Value = 6900
$ws.Range("N97").Value = -7892
$ws.Range("H106").Value = 7042.6665
$ws.Range("I106").Value = 5977.8
$ws.Range("J106").Value = 8373.75
$ws.Range("K106").Value = 5977.8
$ws.Range("L106").Value = 8373.75
$ws.Range("M106").Value = -5346.8
$ws.Range("N106").Value = -9635.75
$ws.Range("H125").Value = 8412.375
$ws.Range("I125").Value = 9000
$ws.Range("J125").Value = 8059.8
$ws.Range("K125").Value = 81000
$ws.Range("L125").Value = 72538.2
$ws.Range("M125").Value = -78540
$ws.Range("N125").Value = -77458.2
$ws.Range("H127").Value = 16628.375
$ws.Range("I127").Value = 12999
$ws.Range("K127").Value = 38997
$ws.Range("M127").Value = -34037
$ws.Range("H132").Value = 415855.47
$ws.Range("I132").Value = 455307.7
$ws.Range("K132").Value = 1365923.1
$ws.Range("M132").Value = -1363393.1
$ws.Range("H135").Value = 6323.3076
$ws.Range("I135").Value = 3020.0667
$ws.Range("J135").Value = 10827.728
$ws.Range("K135").Value = 27180.6003
$ws.Range("L135").Value = 97449.552
$ws.Range("M135").Value = -24645.6003
$ws.Range("N135").Value = -102519.552
$ws.Range("H141").Value = 4739.6
$ws.Range("I141").Value = 3499.5
$ws.Range("J141").Value = 5566.3335
$ws.Range("K141").Value = 10498.5
$ws.Range("L141").Value = 16699.0005
$ws.Range("M141").Value = -5318.5
$ws.Range("N141").Value = -27059.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").Value = $null
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = $null
$ws.Range("H63").Value = 2001
$ws.Range("I63").Value = 2001.2222
$ws.Range("J63").Value = 1999
$ws.Range("K63").Value = 2001.2222
$ws.Range("L63").Value = 1999
$ws.Range("M63").Value = -1315.2222
$ws.Range("N63").Value = -3371
$ws.Range("H66").Value = 2001
$ws.Range("I66").Value = 2001.2222
$ws.Range("J66").Value = 1999
$ws.Range("K66").Value = 10006.111
$ws.Range("L66").Value = 9995
$ws.Range("M66").Value = -6574.110999999999
$ws.Range("N66").Value = -16859
$ws.Range("H122").Value = 2795.3076
$ws.Range("I122").Value = 1941.4
$ws.Range("K122").Value = 5824.200000000001
$ws.Range("M122").Value = -3374.200000000001
$ws.Range("H132").Value = 695840.4
$ws.Range("I132").Value = 734717.9
$ws.Range("K132").Value = 2204153.7
$ws.Range("M132").Value = -2201623.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 900267.4
$ws.Range("I134").Value = 932545.8
$ws.Range("K134").Value = 2797637.4
$ws.Range("M134").Value = -2795102.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6568.1333
$ws.Range("I62").Value = 6794.8
$ws.Range("J62").Value = 6114.8
$ws.Range("K62").Value = 6794.8
$ws.Range("L62").Value = 6114.8
$ws.Range("M62").Value = -6170.8
$ws.Range("N62").Value = -7362.8
$ws.Range("H65").Value = 6568.1333
$ws.Range("I65").Value = 6794.8
$ws.Range("J65").Value = 6114.8
$ws.Range("K65").Value = 33974
$ws.Range("L65").Value = 30574
$ws.Range("M65").Value = -30854
$ws.Range("N65").Value = -36814
$ws.Range("H122").Value = 2231
$ws.Range("I122").Value = 2136.2354
$ws.Range("J122").Value = 2365.25
$ws.Range("K122").Value = 6408.706200000001
$ws.Range("L122").Value = 7095.75
$ws.Range("M122").Value = -3958.706200000001
$ws.Range("N122").Value = -11995.75
$ws.Range("H134").Value = 90917250
$ws.Range("I134").Value = 111117330
$ws.Range("K134").Value = 333351990
$ws.Range("M134").Value = -333349455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 3499.5
$ws.Range("I130").Value = 3499.5
$ws.Range("K130").Value = 10498.5
$ws.Range("M130").Value = -5478.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8189.35
$ws.Range("I7").Value = 8254.75
$ws.Range("K7").Value = 8254.75
$ws.Range("M7").Value = -8142.75
$ws.Range("H40").Value = 4863.8887
$ws.Range("I40").Value = 4734.4
$ws.Range("K40").Value = 4734.4
$ws.Range("M40").Value = -4598.4
$ws.Range("H55").Value = 5073.1763
$ws.Range("I55").Value = 3443.125
$ws.Range("J55").Value = 6522.1113
$ws.Range("K55").Value = 3443.125
$ws.Range("L55").Value = 6522.1113
$ws.Range("M55").Value = -3270.125
$ws.Range("N55").Value = -6868.1113
$ws.Range("H68").Value = 2652.6
$ws.Range("I68").Value = 2039.2858
$ws.Range("J68").Value = 4083.6667
$ws.Range("K68").Value = 2039.2858
$ws.Range("L68").Value = 4083.6667
$ws.Range("M68").Value = -1290.2858
$ws.Range("N68").Value = -5581.6667
$ws.Range("H71").Value = 2652.6
$ws.Range("I71").Value = 2039.2858
$ws.Range("J71").Value = 4083.6667
$ws.Range("K71").Value = 10196.429
$ws.Range("L71").Value = 20418.3335
$ws.Range("M71").Value = -6452.429
$ws.Range("N71").Value = -27906.3335
$ws.Range("H126").Value = 8189.35
$ws.Range("I126").Value = 8254.75
$ws.Range("K126").Value = 24764.25
$ws.Range("M126").Value = -22294.25
$ws.Range("H132").Value = 7835.921
$ws.Range("I132").Value = 5878.3125
$ws.Range("J132").Value = 9259.637000000001
$ws.Range("K132").Value = 17634.9375
$ws.Range("L132").Value = 27778.911
$ws.Range("M132").Value = -15104.9375
$ws.Range("N132").Value = -32838.911

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16650
$ws.Range("I62").Value = 12012.5
$ws.Range("K62").Value = 12012.5
$ws.Range("M62").Value = -11388.5
$ws.Range("H65").Value = 16650
$ws.Range("I65").Value = 12012.5
$ws.Range("K65").Value = 60062.5
$ws.Range("M65").Value = -56942.5
$ws.Range("H96").Value = 1985.625
$ws.Range("I96").Value = 1722
$ws.Range("J96").Value = 2249.25
$ws.Range("K96").Value = 1722
$ws.Range("L96").Value = 2249.25
$ws.Range("M96").Value = -349
$ws.Range("N96").Value = -4995.25
$ws.Range("H113").Value = 9261268
$ws.Range("I113").Value = 13890859
$ws.Range("J113").Value = 2087.5
$ws.Range("K113").Value = 41672577
$ws.Range("L113").Value = 6262.5
$ws.Range("M113").Value = -41670407
$ws.Range("N113").Value = -10602.5
$ws.Range("H126").Value = 2643.3103
$ws.Range("I126").Value = 1551.2632
$ws.Range("K126").Value = 4653.7896
$ws.Range("M126").Value = -2183.7896
